$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 583 ("Femacal de La Calera - Apio" records
# are sorted with the newest entries first). This pushes the former rows 583:616 down to
# 585:618, matching the surrounding data in the sheet.
$ws.Rows("583:584").Insert()

# New row 583: "Primera" quality entry for the new reporting date (2023-04-25 / serial 45041)
$ws.Range("A583").Value = 3
$ws.Range("B583").Value = "Femacal de La Calera"
$ws.Range("C583").Value = "Coquimbo"
$ws.Range("D583").Value = 45041
$ws.Range("E583").Value = 5
$ws.Range("F583").Value = 100112017
$ws.Range("G583").Value = "Apio"
$ws.Range("H583").Value = "Americana (o)"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 120
$ws.Range("K583").Value = 9000
$ws.Range("L583").Value = 9000
$ws.Range("M583").Value = 9000
$ws.Range("N583").Value = "`$/docena de matas"
$ws.Range("O583").Value = "Provincia de Limarí"
$ws.Range("P583").Value = 1500
$ws.Range("Q583").Value = 6
$ws.Range("R583").Value = "Hortaliza"

# New row 584: "Segunda" quality entry for the same new reporting date
$ws.Range("A584").Value = 3
$ws.Range("B584").Value = "Femacal de La Calera"
$ws.Range("C584").Value = "Coquimbo"
$ws.Range("D584").Value = 45041
$ws.Range("E584").Value = 5
$ws.Range("F584").Value = 100112017
$ws.Range("G584").Value = "Apio"
$ws.Range("H584").Value = "Americana (o)"
$ws.Range("I584").Value = "Segunda"
$ws.Range("J584").Value = 110
$ws.Range("K584").Value = 7500
$ws.Range("L584").Value = 7500
$ws.Range("M584").Value = 7500
$ws.Range("N584").Value = "`$/docena de matas"
$ws.Range("O584").Value = "Provincia de Limarí"
$ws.Range("P584").Value = 1250
$ws.Range("Q584").Value = 6
$ws.Range("R584").Value = "Hortaliza"
